# Apply the "New job, old job, P, L, DB columns" edit to the Custom protocol sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row 1 (B1:I1 merged banner) -> extend through column M.
#    Merge first, then re-stamp a single uniform style across the whole
#    range so every cell in the banner keeps the original bold/size-20/
#    centered look (merging alone would let Excel redistribute borders).
# ---------------------------------------------------------------------------
$ws.Range("B1:M1").Merge()
$ws.Range("B1").Copy()
$ws.Range("B1:M1").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Row 2 header labels: new "Job" / "Old Job" columns (style matches the
#    rest of the header row, i.e. C2:J2).
# ---------------------------------------------------------------------------
$ws.Range("J2").Copy()
$ws.Range("L2:M2").PasteSpecial(-4122)
$ws.Range("L2").Value = "Job"
$ws.Range("M2").Value = "Old Job"

# ---------------------------------------------------------------------------
# 3) Data rows 3-90: new "X" / "Y" columns (style matches column D, i.e. s=1).
# ---------------------------------------------------------------------------
$ws.Range("D3").Copy()
$ws.Range("L3:M90").PasteSpecial(-4122)
$ws.Range("L3:L90").Value = "X"
$ws.Range("M3:M90").Value = "Y"

# ---------------------------------------------------------------------------
# 4) Rows 71-90 (Foot Print / Vertical / Lateral / Longitudinal Stiffness
#    "CDTire" block): the inflation-pressure column (D) used to read
#    "IPref" / "0.8\u00b7IPref" / "1.0\u00b7IPref" / "1.2\u00b7IPref" - it is now
#    unified to "P1" for every row in the block.
# ---------------------------------------------------------------------------
$ws.Range("D71:D90").Value = "P1"

# ---------------------------------------------------------------------------
# 5) Header cell K2 ("Displacement [mm]") -> bump its bold font from 11pt
#    to 12pt.
# ---------------------------------------------------------------------------
$ws.Range("K2").Font.Size = 12

# ---------------------------------------------------------------------------
# 6) View state: scrolled down so the frozen pane now starts at row 68 and
#    the active selection sits on F84.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 68
$win.ScrollColumn = 1
$ws.Range("F84").Select()

Write-Host "Edit complete"
